# Apply edits described by the commit "worked on colab and added saved model"
# to the model_study workbook: rename several row/column labels to their
# English, cleaned-up equivalents and add a new "Test 6" column (G) with
# results for an additional experiment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing labels in column A ---
$ws.Range("A2").Value  = "Balanced dataset (undersampled labels)"
$ws.Range("A5").Value  = "Autofeat generated Features"
$ws.Range("A6").Value  = "Other Features"
$ws.Range("A8").Value  = "Other feautres Normalized"
$ws.Range("A9").Value  = "Number of layers"
$ws.Range("A11").Value = "Optimizer"
$ws.Range("A12").Value = "Epochs"
$ws.Range("A13").Value = "Batch size"
$ws.Range("A15").Value = "Results"
$ws.Range("A16").Value = "Training time"
$ws.Range("A17").Value = "It overfitted"
$ws.Range("A18").Value = "Accuracy training"
$ws.Range("A19").Value = "Loss training"
$ws.Range("A20").Value = "Accuracy validation"
$ws.Range("A21").Value = "Loss validation"
$ws.Range("A22").Value = "Accuracy test"
$ws.Range("A23").Value = "Loss Test"
$ws.Range("A24").Value = "Confussion Matrix"

# --- Add new "Test 6" column G with its results ---
# Header cell: bold like the rest of row 1, general (non-centered) alignment.
$ws.Range("G1").Value = "Test 6"
$ws.Range("G1").Font.Bold = $true

# Data cells: centered both horizontally and vertically, like column F.
$ws.Range("G2").Value = "Yes"
$ws.Range("G3").Value = "Yes"
$ws.Range("G4").Value = "Yes"
$ws.Range("G5").Value = "Yes"
$ws.Range("G6").Value = "Yes (**2, **3, exp)"
$ws.Range("G7").Value = "Yes"
$ws.Range("G8").Value = "Yes"

$ws.Range("G2:G3").HorizontalAlignment = -4108
$ws.Range("G2:G3").VerticalAlignment = -4108
$ws.Range("G5:G8").HorizontalAlignment = -4108
$ws.Range("G5:G8").VerticalAlignment = -4108

# G4 only got a plain centered alignment (no vertical centering) in the
# original edit, matching the rest of row 4's mixed formatting.
$ws.Range("G4").HorizontalAlignment = -4108

# --- Column widths (best-fit, as Excel would do after widening for the
# longer labels / new column) ---
$ws.Columns.Item(1).ColumnWidth = 37.28515625
$ws.Columns.Item(7).ColumnWidth = 17.42578125

# --- Selection, matching the saved view state ---
$ws.Range("G9").Select() | Out-Null
